$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date header in B1, keep it stored as text (not auto-converted to a date serial)
$ws.Range("B1").Value = "18/03/2023"

# Update the hourly values in column B (rows 2-17)
$ws.Range("B2").Value = 73
$ws.Range("B3").Value = 153
$ws.Range("B4").Value = 128
$ws.Range("B5").Value = 117
$ws.Range("B6").Value = 109
$ws.Range("B7").Value = 93
$ws.Range("B8").Value = 91
$ws.Range("B9").Value = 118
$ws.Range("B10").Value = 97
$ws.Range("B11").Value = 99
$ws.Range("B12").Value = 85
$ws.Range("B13").Value = 76
$ws.Range("B14").Value = 61
$ws.Range("B15").Value = 24
$ws.Range("B16").Value = 24
$ws.Range("B17").Value = 19
